$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 79; this shifts the existing rows 79..192
# down to 80..193 (carrying their formatting/values with them), matching
# the diff's observed "every record from row 79 on moves down by one, a
# brand-new record appears at row 79" shape.
$ws.Rows(79).Insert()

# Populate the newly inserted row 79 with the new weekly data point.
$ws.Cells.Item(79, 1).Value = 4
$ws.Cells.Item(79, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(79, 3).Value = "Los Lagos"
$ws.Cells.Item(79, 4).Value = 44579
$ws.Cells.Item(79, 4).NumberFormat = $ws.Cells.Item(80, 4).NumberFormat
$ws.Cells.Item(79, 5).Value = 10
$ws.Cells.Item(79, 6).Value = "Fruta"
$ws.Cells.Item(79, 7).Value = 100108
$ws.Cells.Item(79, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(79, 9).Value = 100108005
$ws.Cells.Item(79, 10).Value = "Piña"
$ws.Cells.Item(79, 11).Value = "Caramelo"
$ws.Cells.Item(79, 12).Value = "Tercera"
$ws.Cells.Item(79, 13).Value = 200
$ws.Cells.Item(79, 14).Value = 19000
$ws.Cells.Item(79, 15).Value = 20000
$ws.Cells.Item(79, 16).Value = 19500
$ws.Cells.Item(79, 17).Value = "`$/caja 16 unidades"
$ws.Cells.Item(79, 18).Value = "Ecuador"
$ws.Cells.Item(79, 19).Value = 1219
$ws.Cells.Item(79, 20).Value = 16
